$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we touch so numeric-looking strings
# like "1.000" or "5.399" are preserved verbatim instead of being
# auto-converted to numbers by Excel's smart-entry parsing.
$targetRefs = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","D8","E8","D9","E9","E10","E11","D12","E12","D13","E13","D14","E14","E15","D16","E16","D17","E17","D18","B19","C19","D19","E19","B20","C20","D20","E20","D21","E21","D22","E22","E23","D24","E24","D25","E25","D26","E26","E28","D29","E29","D30","E30","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","B46","C46","D46","E46","B47","C47","D47","E47","B48","C48","D48","E48","B49","C49","D49","E49","D50","E50","D51","E51")
foreach ($ref in $targetRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated coin data (prices, 1h volume deltas, and the
# re-ranked rows 19/20 and 46-49) to match the latest GitHub Actions
# refresh of the cryptos list.
$ws.Range("D2").Value = "29.899.53"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.893.15"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "0.7720"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "244.12"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.3125"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").Value = "25.65"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  +7.53%  "
$ws.Range("D12").Value = "1.999.87"
$ws.Range("E12").Value = "  +4.68%  "
$ws.Range("D13").Value = "0.7703"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "5.399"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "6.214"
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("D17").Value = "30.161.78"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "13.93"
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.327.47"
$ws.Range("E19").Value = "  +6.27%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "245.16"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "0.000007854"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").Value = "8.179"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "0.1593"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("D26").Value = "9.505"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").Value = "2.040"
$ws.Range("E29").Value = "  -0.69%  "
$ws.Range("D30").Value = "1.439"
$ws.Range("E30").Value = "  +2.54%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "4.520"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "4.114"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "0.05468"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "1.248"
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("D36").Value = "0.7544"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "1.004"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "2.711"
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("D39").Value = "0.01968"
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "0.4504"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "73.85"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "1.094.73"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "6.048"
$ws.Range("E44").Value = "  +2.83%  "
$ws.Range("D45").Value = "0.8548"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.202.32"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "103.03"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.882"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "7.629"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").Value = "9.820"
$ws.Range("E51").Value = "  -2.07%  "
